# Update the cryptocurrency price ("Price", col D) and volume change
# ("Volume(1h)", col E) figures in-place, exactly as produced by the
# latest scrape. Values that look like plain decimal numbers (e.g.
# "1.004") must be forced to remain TEXT (matching the original
# inlineStr/shared-string cells) instead of being auto-converted to a
# number by Excel, which would also silently drop formatting such as
# trailing zeros. We do this by temporarily marking the cell as Text
# ("@") before assigning the value, then resetting the cell style back
# to "Normal" afterwards so no stray per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "22.414.79"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.572.22"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3735"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3391"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07549"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.134"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.989"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.922"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "1.571.36"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06737"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.274"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "22.426.67"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.628"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.007"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "1.749.20"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.380"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02460"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06509"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.448"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6208"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.813"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5791"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.073"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.215"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07317"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
